# Update "想去人数" (wanted attendance) figures on the 展览 (Exhibition) sheet
# and the 全部类型 (All types) aggregate sheet, per the generated-output refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - rows 2-7 column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 14720
$wsExpo.Range("F3").Value = 341
$wsExpo.Range("F4").Value = 703
$wsExpo.Range("F5").Value = 243
$wsExpo.Range("F6").Value = 610
$wsExpo.Range("F7").Value = 1570

# Sheet "全部类型" (all types, aggregated) - matching rows (2-5, 8-9) column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14720
$wsAll.Range("F3").Value = 341
$wsAll.Range("F4").Value = 703
$wsAll.Range("F5").Value = 243
$wsAll.Range("F8").Value = 610
$wsAll.Range("F9").Value = 1570
